$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45174 = 2023-09-05) for
# each data row (rows 2-21). Bump it forward by one day (to serial 45175 =
# 2023-09-06) for every row, preserving the existing date formatting.
for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
